$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.455.07'
$ws.Range("E2").Value = '  +0.91%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.787.60'
$ws.Range("E3").Value = '  +0.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '697.80'
$ws.Range("E5").Value = '  +6.58%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.46'
$ws.Range("E6").Value = '  +3.15%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.785.91'
$ws.Range("E7").Value = '  +0.17%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  -0.28%  '

# Row 10
$ws.Range("E10").Value = '  +1.36%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.29'
$ws.Range("E11").Value = '  +5.13%  '

# Row 12
$ws.Range("E12").Value = '  -0.01%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  +6.10%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.99'
$ws.Range("E14").Value = '  +2.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.429.61'
$ws.Range("E15").Value = '  +0.13%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.794.46'
$ws.Range("E16").Value = '  +0.36%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.491.58'
$ws.Range("E17").Value = '  +0.96%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.65'
$ws.Range("E18").Value = '  -0.36%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.16'
$ws.Range("E19").Value = '  +1.59%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '478.14'
$ws.Range("E22").Value = '  +1.49%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.709'
$ws.Range("E23").Value = '  -0.22%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.66'
$ws.Range("E24").Value = '  +1.66%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000141'
$ws.Range("E25").Value = '  -1.94%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.27'
$ws.Range("E26").Value = '  -0.53%  '

# Row 27
$ws.Range("E27").Value = '  +0.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.14'
$ws.Range("E28").Value = '  +0.76%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.942.62'
$ws.Range("E29").Value = '  +0.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.13%  '

# Row 31
$ws.Range("E31").Value = '  +13.76%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.50'
$ws.Range("E32").Value = '  +4.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.27'
$ws.Range("E33").Value = '  -0.18%  '

# Row 34
$ws.Range("E34").Value = '  +8.15%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.25'
$ws.Range("E35").Value = '  +1.38%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.21'
$ws.Range("E36").Value = '  +3.37%  '

# Row 37
$ws.Range("E37").Value = '  +0.24%  '

# Row 38
$ws.Range("E38").Value = '  +1.36%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("E39").Value = '  +1.93%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.97'
$ws.Range("E40").Value = '  +1.88%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("E41").Value = '  +9.92%  '

# Row 42
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.976'
$ws.Range("E42").Value = '  +1.78%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.18%  '

# Row 44
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("B45").Value = 'FLOKI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000320'
$ws.Range("E45").Value = '  +18.27%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '163.30'
$ws.Range("E46").Value = '  +3.88%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.77'
$ws.Range("E47").Value = '  +1.98%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.42'
$ws.Range("E48").Value = '  -3.39%  '

# Row 49
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.38'
$ws.Range("E49").Value = '  -0.93%  '

# Row 50
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.299'
$ws.Range("E50").Value = '  -0.35%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.56'
$ws.Range("E51").Value = '  +1.53%  '
